$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:B73 (rows 2 and 3 unchanged, 4-73 updated per diff)
$values = @{
    2 = "0.3497567887281834"
    3 = "-0.06966572625950762"
    4 = "-0.5718192715580426"
    5 = "0.02513501597327435"
    6 = "-0.6777052466944516"
    7 = "-0.04687708722170142"
    8 = "-0.6491824449565908"
    9 = "-0.03143360567770964"
    10 = "0.3822062224503651"
    11 = "-0.04722993886919886"
    12 = "-0.5591348797726128"
    13 = "-0.3824674568541571"
    14 = "0.3"
    15 = "0"
    16 = "0"
    17 = "-0.07670146240877694"
    18 = "-0.2"
    19 = "0.1"
    20 = "0"
    21 = "-0.1"
    22 = "0.1"
    23 = "0"
    24 = "0.3"
    25 = "-0.3"
    26 = "-0.2"
    27 = "0"
    28 = "0.07231168587756959"
    29 = "-0.06367835427788604"
    30 = "0.08504296140338079"
    31 = "0.0003037193523987569"
    32 = "-0.01801114345337739"
    33 = "0.03826442504165321"
    34 = "-0.02232391446496779"
    35 = "-0.06208406292622634"
    36 = "0.006968764489216402"
    37 = "-0.121948034731228"
    38 = "-0.0361595126957486"
    39 = "-0.08573317974146363"
    40 = "0.004670068289308601"
    41 = "0.04099992594395328"
    42 = "0.05927521680271484"
    43 = "0.04061467511413108"
    44 = "0.08363338226171732"
    45 = "-0.184542838621156"
    46 = "-0.4"
    47 = "-0.3200594507515429"
    48 = "0.1"
    49 = "-0.2628511487790233"
    50 = "-0.3133720122661878"
    51 = "0.5469624758293199"
    52 = "1.1"
    53 = "-0.4440571223929872"
    54 = "-0.7255945204468831"
    55 = "-0.5292660609007143"
    56 = "-0.1550786956675604"
    57 = "-2.168330733759602"
    58 = "-0.03982694963614287"
    59 = "0.2669401745841223"
    60 = "0.03791487406588956"
    61 = "-0.04567208272808071"
    62 = "-0.5154625125417773"
    63 = "-0.1813602613933202"
    64 = "-0.01480819732384536"
    65 = "0.02918400950819283"
    66 = "-0.03321544329283629"
    67 = "1.303303454188581E-05"
    68 = "-0.006125572440376981"
    69 = "0.04879937325030748"
    70 = "0.0477695913607396"
    71 = "0.4714513528429705"
    72 = "-0.02605454389395597"
    73 = "0.04549112474043772"
}

foreach ($key in $values.Keys) {
    $ws.Cells.Item([int]$key, 2).Value = [double]$values[$key]
}

# Remove rows 74-82 (data series shortened)
$ws.Rows("74:82").Delete()
